$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210, shifting existing rows 210:275 down to 211:276
$ws.Rows("210:210").Insert()

# Populate the new row 210 with the new weekly price observation
$ws.Range("A210").Value = 4
$ws.Range("B210").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C210").Value = "Los Lagos"
$ws.Range("D210").Value = 44588
$ws.Range("E210").Value = 10
$ws.Range("F210").Value = 100112023
$ws.Range("G210").Value = "Brócoli"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 300
$ws.Range("K210").Value = 1500
$ws.Range("L210").Value = 1500
$ws.Range("M210").Value = 1500
$ws.Range("N210").Value = "$/unidad"
$ws.Range("O210").Value = "Región Metropolitana"
$ws.Range("P210").Value = 1500
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = "Hortaliza"
